# Mail-showcase workbook edit:
# Adds a new "localdb" command group to the hidden '#system' sheet that backs
# the named ranges used for auto-complete/validation of script "target" /
# "command" values. This mirrors inserting a new column (N) on the
# '#system' sheet, populating it with the localdb function list, inserting
# "localdb" into the alphabetical category list in column A, and fixing up
# every defined name whose range shifted as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Insert a brand-new column at N; everything from N..AC slides right to
#    O..AD (all 127 rows), exactly like Excel's own "Insert Column" command.
$ws.Columns("N:N").Insert()

# 2) Insert the new "localdb" entry into the alphabetically sorted category
#    list in column A, right before "macro" (row 14). Only column A's
#    values shift down by one row; every other column on the sheet is
#    untouched, so the cells are rewritten directly rather than using
#    Range.Insert (which would shift entire row bands on this engine).
$ws.Range("A30").Value = "xml"
$ws.Range("A29").Value = "ws.async"
$ws.Range("A28").Value = "ws"
$ws.Range("A27").Value = "webcookie"
$ws.Range("A26").Value = "webalert"
$ws.Range("A25").Value = "web"
$ws.Range("A24").Value = "step"
$ws.Range("A23").Value = "ssh"
$ws.Range("A22").Value = "sound"
$ws.Range("A21").Value = "sms"
$ws.Range("A20").Value = "redis"
$ws.Range("A19").Value = "rdbms"
$ws.Range("A18").Value = "pdf"
$ws.Range("A17").Value = "number"
$ws.Range("A16").Value = "mail"
$ws.Range("A15").Value = "macro"
$ws.Range("A14").Value = "localdb"

# 3) Populate the new column N with the localdb header + its function list.
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# 4) Fix up the defined names whose referenced column moved one letter to
#    the right because of the column insert at N, and extend "target" by
#    the one extra row added to column A.
$wb.Names.Item("mail").RefersTo       = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo     = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo        = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo      = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo      = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("ssh").RefersTo        = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo       = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo     = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo        = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("sound").RefersTo      = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("sms").RefersTo        = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("macro").RefersTo      = "='#system'!`$O`$2:`$O`$4"

# 5) Register the new "localdb" named range.
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
